# Apply updated employee absence data values to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 65455
$ws.Range("B2").Value = "Isabelly Monteiro"
$ws.Range("C2").Value = "P&D"
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 45101
$ws.Range("G2").Value = 5434.74

# Row 3
$ws.Range("A3").Value = 3995
$ws.Range("B3").Value = "Lunna Gomes"
$ws.Range("C3").Value = "P&D"
$ws.Range("D3").Value = "Doenca"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 45079
$ws.Range("G3").Value = 2291.97

# Row 4
$ws.Range("A4").Value = 12661
$ws.Range("B4").Value = "Igor Cardoso"
$ws.Range("C4").Value = "Operacoes"
$ws.Range("D4").Value = "Consulta medica"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 45084
$ws.Range("G4").Value = 4947.65

# Row 5
$ws.Range("A5").Value = 1953
$ws.Range("B5").Value = "Marcos Vinicius Araújo"
$ws.Range("C5").Value = "Operacoes"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 45089
$ws.Range("G5").Value = 9394.629999999999

# Row 6
$ws.Range("A6").Value = 78962
$ws.Range("B6").Value = "Maria Isis Pereira"
$ws.Range("D6").Value = "Consulta medica"
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 45105
$ws.Range("G6").Value = 7159.66

# Row 7
$ws.Range("A7").Value = 49747
$ws.Range("B7").Value = "João Guilherme Monteiro"
$ws.Range("C7").Value = "Atendimento ao Cliente"
$ws.Range("D7").Value = "Outros"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 45096
$ws.Range("G7").Value = 3262.44

# Row 8
$ws.Range("A8").Value = 60266
$ws.Range("B8").Value = "Marcela Dias"
$ws.Range("C8").Value = "Operacoes"
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 45084
$ws.Range("G8").Value = 7508.92

# Row 9
$ws.Range("A9").Value = 60615
$ws.Range("B9").Value = "Sr. Luiz Miguel Cirino"
$ws.Range("C9").Value = "Marketing"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 45094
$ws.Range("G9").Value = 5809.08

# Row 10
$ws.Range("A10").Value = 70295
$ws.Range("B10").Value = "Manuella Garcia"
$ws.Range("C10").Value = "Juridico"
$ws.Range("D10").Value = "Consulta medica"
$ws.Range("F10").Value = 45092
$ws.Range("G10").Value = 4106.68

# Row 11
$ws.Range("A11").Value = 28069
$ws.Range("B11").Value = "Olívia Martins"
$ws.Range("D11").Value = "Outros"
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 45095
$ws.Range("G11").Value = 4555.22
